# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "K" (strikeouts) values replacing the old "Strike#" values in column G
# for rows 2-28 (data rows). Header in G1 stays "K".
$kValues = @(6, 2, 4, 7, 5, 5, 5, 2, 3, 6, 6, 4, 8, 5, 4, 5, 2, 7, 6, 4, 6, 4, 6, 5, 4, 4, 2)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
